$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: convert text dates to real date values (date-formatted) ---
# Rows 2-7 currently hold text dates ("2025-09-15"/"2025-09-16") with no
# explicit number format; give them numeric serials and the existing date
# style (copied from A8, which is already date-formatted) so we reuse the
# workbook's existing style index instead of inventing a new one.
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A2:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A2").Value = 45914
$ws.Range("A3").Value = 45914
$ws.Range("A4").Value = 45915
$ws.Range("A5").Value = 45916
$ws.Range("A6").Value = 45916
$ws.Range("A7").Value = 45916

# Rows 8-13 already use the date style; just correct the serial values.
$ws.Range("A8").Value = 45917
$ws.Range("A9").Value = 45917
$ws.Range("A10").Value = 45917
$ws.Range("A11").Value = 45918
$ws.Range("A12").Value = 45918
$ws.Range("A13").Value = 45918

# --- New food/portion entries, written in the same order the new unique
# strings first appear so the shared-string table matches. ---

# Row 10: extra meal columns
$ws.Range("G10").Value = 300
$ws.Range("H10").Value = "zsömle"

# D3: rename food item "csirke mell" -> "csirke"
$ws.Range("D3").Value = "csirke"

# Row 11: extra meal columns
$ws.Range("E11").Value = 300
$ws.Range("F11").Value = "csirke"
$ws.Range("G11").Value = 400
$ws.Range("H11").Value = "krumpli"

# Row 12: extra meal columns
$ws.Range("E12").Value = 400
$ws.Range("F12").Value = "hal"
$ws.Range("G12").Value = 300
$ws.Range("H12").Value = "rizs"

# Row 13: extra meal columns
$ws.Range("E13").Value = 300
$ws.Range("F13").Value = "sertés"
$ws.Range("G13").Value = 400
$ws.Range("H13").Value = "krumpli"

# Row 8: extra meal columns (all strings already exist in the table)
$ws.Range("G8").Value = 400
$ws.Range("H8").Value = "ponty"
$ws.Range("I8").Value = 400
$ws.Range("J8").Value = "kenyér"

# --- Selection moves to D3 ---
$ws.Range("D3").Select() | Out-Null
